# Removed components that only had one file, made into objects.
#
# For each "Object" row that has exactly one following "Component" row
# (sharing the same Object Unique ID in column A), copy the component's
# File name / File use / Type of Resource / Language (columns C:F) up
# into the Object row, then delete the now-redundant Component row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Item description")

# Rows (in the original layout) that are single-component objects and
# therefore get merged: the Object row immediately followed by its only
# Component row.
$pairs = @(
    @{Obj = 6;  Comp = 7},
    @{Obj = 8;  Comp = 9},
    @{Obj = 10; Comp = 11},
    @{Obj = 12; Comp = 13},
    @{Obj = 14; Comp = 15},
    @{Obj = 16; Comp = 17},
    @{Obj = 18; Comp = 19},
    @{Obj = 20; Comp = 21},
    @{Obj = 22; Comp = 23},
    @{Obj = 24; Comp = 25}
)

foreach ($pair in $pairs) {
    $objRow = $pair.Obj
    $compRow = $pair.Comp

    $srcRange = $ws.Range("C$compRow`:F$compRow")
    $dstRange = $ws.Range("C$objRow`:F$objRow")

    $srcRange.Copy()
    $dstRange.PasteSpecial(-4104)  # xlPasteAll
}

$ws.Application.CutCopyMode = $false

# Delete the component rows, from the bottom up so row numbers of rows
# still to be processed don't shift.
$compRows = $pairs | ForEach-Object { $_.Comp } | Sort-Object -Descending
foreach ($r in $compRows) {
    $ws.Rows.Item($r).Delete()
}

# Reset the selection as recorded in the saved workbook.
$ws.Range("C2").Select()
